# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old rows that are no longer used in the new layout ---
# (old "sector distribution / data block" that lived in rows 5-9, plus the
# old row 13 "Sector Distribution Details" label, all get repositioned
# below, so clear them out first)
$ws.Range("A5:D9").Clear()
$ws.Range("A13:D13").Clear()

# --- Helper-ish inline writes for the new layout ---

# Row 1 - name (unchanged)
$c = $ws.Range("A1")
$c.Value = "Bosnia and Herzegovina"
$c.Font.Size = 18
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false

# Row 3 - title (unchanged)
$c = $ws.Range("A3")
$c.Value = "MSME Participation on the Economy"
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false

# Row 8 - NEW bold+underline "Source Type" header
$c = $ws.Range("A8")
$c.Value = "Source Type: Ministry of Finance/Central Bank (Most Widely Used)"
$c.Font.Bold = $true
$c.Font.Underline = $true

# Row 10 - column headers (Micro / SMEs / MSMEs)
$c = $ws.Range("B10"); $c.Value = "Micro"; $c.Font.Bold = $true
$c = $ws.Range("C10"); $c.Value = "SMEs";  $c.Font.Bold = $true
$c = $ws.Range("D10"); $c.Value = "MSMEs"; $c.Font.Bold = $true

# Row 11 - Enterprises (absolute #)
$c = $ws.Range("A11"); $c.Value = "Enterprises (absolute #)"; $c.Font.Bold = $true
$c = $ws.Range("B11"); $c.NumberFormat = "@"; $c.Value = "151107"
$c = $ws.Range("C11"); $c.NumberFormat = "@"; $c.Value = "10188"
$c = $ws.Range("D11"); $c.NumberFormat = "@"; $c.Value = "161295"

# Row 12 - Enterprises density (per 1000 people)
$c = $ws.Range("A12"); $c.Value = "Enterprises density (per 1000 people)"; $c.Font.Bold = $true
$c = $ws.Range("B12"); $c.NumberFormat = "@"; $c.Value = "39.1"
$c = $ws.Range("C12"); $c.NumberFormat = "@"; $c.Value = "2.6"
$c = $ws.Range("D12"); $c.NumberFormat = "@"; $c.Value = "41.8"

# Row 13 - NEW Employment (% of total)
$c = $ws.Range("A13"); $c.Value = "Employment (% of total)"; $c.Font.Bold = $true
$c = $ws.Range("D13"); $c.NumberFormat = "@"; $c.Value = "48.8"

# Row 14 - Enterprises (% of total)
$c = $ws.Range("A14"); $c.Value = "Enterprises (% of total)"; $c.Font.Bold = $true
$c = $ws.Range("B14"); $c.NumberFormat = "@"; $c.Value = "93.3"
$c = $ws.Range("C14"); $c.NumberFormat = "@"; $c.Value = "6.3"
$c = $ws.Range("D14"); $c.NumberFormat = "@"; $c.Value = "99.6"

# Row 15 - source note (italic) - moved from old row 9
$c = $ws.Range("A15")
$c.Value = "Source: MVTEO, 2008"
$c.Font.Italic = $true
$c.Font.Bold = $false
$c.Font.Underline = $false

# Row 20 - Sector Distribution Details (moved from old row 13)
$c = $ws.Range("A20")
$c.Value = "Sector Distribution Details"
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false

# Row 23 - NEW "MVTEO" title
$c = $ws.Range("A23")
$c.Value = "MVTEO"
$c.Font.Bold = $true

# Row 24 - NEW citation (italic, "source" style)
$c = $ws.Range("A24")
$c.Value = 'Ministry of Foreign Trade and Economic Relations of Bosnia and Herzegovina (MVTEO), "Small and Medium-Sized Enterprise Development Strategy in Bosnia and Herzegovina 2009 - 2011", 2009. Available at http://www.mvteo.gov.ba/vijesti/posljednje_vijesti/default.aspx?id=1204&langTag=bs-BA'
$c.Font.Italic = $true
$c.Font.Bold = $false
$c.Font.Underline = $false

# --- Rename the sheet from "Data" to "Summary" ---
$ws.Name = "Summary"

Write-Host "Edit complete"
